$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to remain plain text (values like "19.45" would
# otherwise be auto-converted to numbers by Excel's type inference), matching
# the inlineStr/shared-string representation used throughout the sheet.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '44.133.30'
$ws.Range("E2").Value = '  +1.49%  '
$ws.Range("D3").Value = '2.353.43'
$ws.Range("E3").Value = '  -0.89%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '0.678'
$ws.Range("E5").Value = '  +5.23%  '
$ws.Range("D6").Value = '241.17'
$ws.Range("E6").Value = '  +3.05%  '
$ws.Range("D7").Value = '74.01'
$ws.Range("E7").Value = '  +6.65%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").Value = '0.564'
$ws.Range("E9").Value = '  +23.00%  '
$ws.Range("D10").Value = '0.102'
$ws.Range("E10").Value = '  +5.41%  '
$ws.Range("D11").Value = '31.35'
$ws.Range("E11").Value = '  +18.58%  '
$ws.Range("D12").Value = '7.42'
$ws.Range("E12").Value = '  +19.89%  '
$ws.Range("E13").Value = '  +1.96%  '
$ws.Range("D14").Value = '2.702.94'
$ws.Range("E14").Value = '  -0.90%  '
$ws.Range("D15").Value = '16.83'
$ws.Range("E15").Value = '  +7.42%  '
$ws.Range("D16").Value = '0.910'
$ws.Range("E16").Value = '  +6.83%  '
$ws.Range("D17").Value = '2.350.14'
$ws.Range("E17").Value = '  -1.03%  '
$ws.Range("D18").Value = '44.381.88'
$ws.Range("E18").Value = '  +2.09%  '
$ws.Range("D19").Value = '0.0000103'
$ws.Range("E19").Value = '  +3.88%  '
$ws.Range("D20").Value = '6.69'
$ws.Range("E20").Value = '  +5.78%  '
$ws.Range("D21").Value = '77.93'
$ws.Range("E21").Value = '  +5.34%  '
$ws.Range("D22").Value = '255.96'
$ws.Range("E22").Value = '  +2.03%  '
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").Value = '3.78'
$ws.Range("E24").Value = '  -4.82%  '
$ws.Range("D25").Value = '2.55'
$ws.Range("E25").Value = '  +4.09%  '
$ws.Range("E26").Value = '  +7.36%  '
$ws.Range("E27").Value = '  +1.38%  '
$ws.Range("D28").Value = '22.57'
$ws.Range("E28").Value = '  -1.22%  '
$ws.Range("D29").Value = '174.53'
$ws.Range("E29").Value = '  +1.12%  '
$ws.Range("D30").Value = '1.59'
$ws.Range("E30").Value = '  +2.83%  '
$ws.Range("D31").Value = '0.131'
$ws.Range("E31").Value = '  +3.47%  '
$ws.Range("E32").Value = '  +4.67%  '
$ws.Range("D33").Value = '5.38'
$ws.Range("E33").Value = '  +7.93%  '
$ws.Range("D34").Value = '0.0750'
$ws.Range("E34").Value = '  +8.83%  '
$ws.Range("D35").Value = '5.35'
$ws.Range("E35").Value = '  +4.97%  '
$ws.Range("D36").Value = '3.90'
$ws.Range("E36").Value = '  +7.62%  '
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("D38").Value = '6.54'
$ws.Range("E38").Value = '  -0.71%  '
$ws.Range("D39").Value = '0.0274'
$ws.Range("E39").Value = '  +7.64%  '
$ws.Range("D40").Value = '19.45'
$ws.Range("E40").Value = '  +5.43%  '
$ws.Range("D41").Value = '8.98'
$ws.Range("E41").Value = '  +0.06%  '
$ws.Range("E42").Value = '  -0.08%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '1.26'
$ws.Range("E43").Value = '  +3.51%  '
$ws.Range("B44").Value = 'Cronos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D44").Value = '0.0998'
$ws.Range("E44").Value = '  +4.80%  '
$ws.Range("D45").Value = '0.189'
$ws.Range("E45").Value = '  +13.43%  '
$ws.Range("D46").Value = '100.16'
$ws.Range("E46").Value = '  +1.08%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '2.45'
$ws.Range("E47").Value = '  +9.35%  '
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").Value = '1.17'
$ws.Range("E48").Value = '  -1.83%  '
$ws.Range("E49").Value = '  -0.10%  '
$ws.Range("D50").Value = '1.451.90'
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("D51").Value = '2.80'
$ws.Range("E51").Value = '  +2.18%  '

# Restore the default (Normal) style on the price column so no stray number
# formatting is left applied to the cells themselves.
$priceRange.Style = "Normal"

